# Custom properties - bug fixes
$wb = $excel.ActiveWorkbook

# --- AMSIN sheet: append a new registration-history row (row 15) ---
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Cells.Item(15, 1).Value = "2023-04-18"
$wsAmsin.Cells.Item(15, 2).Value = 45034.50914052363
$wsAmsin.Cells.Item(15, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAmsin.Cells.Item(15, 3).Value = "176fstrrun"
$wsAmsin.Cells.Item(15, 4).Value = 46
$wsAmsin.Cells.Item(15, 5).Value = 46
$wsAmsin.Cells.Item(15, 6).Value = 0
$wsAmsin.Cells.Item(15, 7).Value = 0.6899999999999999

# --- AMS sheet: fix up row 12 formatting / precise run time ---
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(12, 2).Value = 45027.79182425926

$wsAms.Range("A12").Style = "Normal"
$wsAms.Range("C12:G12").Style = "Normal"
